# DsLookup.xlsx refactor (#67, #71): wrap the DSLOOKUP() calls with INDEX()
# and add a second example row demonstrating INDEX(...,2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F7 used to read:  =DSLOOKUP("Sheet1", "Family", "Smith", "Age" )
# it now reads:     =INDEX(DSLOOKUP("Sheet1", "Family", "Smith", "Age" ),1)
$ws.Range("F7").Formula = '=INDEX(DSLOOKUP("Sheet1", "Family", "Smith", "Age" ),1)'

# New row 8 / F8 demonstrates the second element of the (still erroring)
# DSLOOKUP result: =INDEX(DSLOOKUP("Sheet1", "Family", "Smith", "Age" ),2)
$ws.Range("F8").Formula = '=INDEX(DSLOOKUP("Sheet1", "Family", "Smith", "Age" ),2)'

# The sheet's dimension (A1:F7 -> A1:F8) is recomputed automatically by the
# engine once F8 has content.

# Move the active selection to the next empty cell below the new data (F9),
# matching the author's saved cursor position.
$ws.Range("F9").Select()

# The workbook window had also been scrolled/moved on screen
# (yWindow 105 -> 2805) before saving; reflect that on the automation
# window object as well.
$excel.ActiveWindow.Top = 2805
